$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 used to hold the first student's data; it becomes a second
# "header" row (same 12 labels as row 1, same row height / style).
# ClearContents() first so every cell in the row counts as brand new -
# new cells pick up the column's default style (index 2, the same
# style row 1's header cells use), matching the diff's s="2" on A2:L2.
# ---------------------------------------------------------------------
$ws.Range("A2:L2").ClearContents()
$headers = @("name", "msv", "class", "hdcm.uv1", "hdcm.uv2", "hdcm.uv3", `
             "hdcm.uv4", "hdcm.uv5", "hd.01", "hd.02", "hd.03", "pb")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(2, $col).Value = $headers[$col - 1]
}
$ws.Rows.Item(2).RowHeight = 18.75

# ---------------------------------------------------------------------
# Row 3 used to hold the second student's data; it is replaced with a
# new student's row. A3:D3 already exist, so writing .Value in place
# keeps them free of any explicit cell style (matching the diff, which
# leaves A3:D3 without an s="..." attribute). E3:L3 are cleared to
# blank/empty - new cells auto-inherit the column style, so each is
# reset back to the "Normal" style right after, keeping them styleless
# like the diff.
# ---------------------------------------------------------------------
$ws.Cells.Item(3, 1).Value = "Đỗ Trọng Khôi"
$ws.Cells.Item(3, 2).Value = "B20DCDT112"
$ws.Cells.Item(3, 3).Value = "D20DTMT2"
$ws.Cells.Item(3, 4).Value = "Trần Tuấn Anh - C3.3: 3 - C4.2: 4 - C5.3: 5 - C6.3: 6 - C6.4: 7 - GPA: 8"

for ($col = 5; $col -le 12; $col++) {
    $cell = $ws.Cells.Item(3, $col)
    $cell.Value = ""
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Table1 (and its AutoFilter) grows by one row, now spanning the new
# second header-like row too: A1:L1 -> A1:L2.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:L2"))

Write-Output "done"
